$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$items = @(
    @("Double Cheese Burger", 2),
    @("Chicken Burger", 1),
    @("MEGA BURGER", 1),
    @("Cheese Burger", 1)
)

$timestamp = "2024-12-10 03:28:44"

$row = 6
foreach ($item in $items) {
    $ws.Cells.Item($row, 1).Value = $item[0]
    $ws.Cells.Item($row, 2).Value = $item[1]
    $ws.Cells.Item($row, 3).Value = $timestamp
    $row = $row + 1
}
